$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "72.256.86"
Set-TextValue $ws.Range("E2") "  +4.21%  "

Set-TextValue $ws.Range("D3") "3.615.70"
Set-TextValue $ws.Range("E3") "  +6.74%  "

Set-TextValue $ws.Range("D4") "0.998"
Set-TextValue $ws.Range("E4") "  -0.19%  "

Set-TextValue $ws.Range("D5") "590.14"
Set-TextValue $ws.Range("E5") "  +0.38%  "

Set-TextValue $ws.Range("D6") "184.13"
Set-TextValue $ws.Range("E6") "  +2.48%  "

Set-TextValue $ws.Range("D7") "3.605.12"
Set-TextValue $ws.Range("E7") "  +6.58%  "

Set-TextValue $ws.Range("D8") "0.608"
Set-TextValue $ws.Range("E8") "  +1.84%  "

Set-TextValue $ws.Range("D9") "0.999"
Set-TextValue $ws.Range("E9") "  -0.04%  "

Set-TextValue $ws.Range("E10") "  +6.84%  "

Set-TextValue $ws.Range("D11") "0.610"
Set-TextValue $ws.Range("E11") "  +3.38%  "

Set-TextValue $ws.Range("D12") "50.25"
Set-TextValue $ws.Range("E12") "  +3.79%  "

Set-TextValue $ws.Range("D13") "0.0000293"
Set-TextValue $ws.Range("E13") "  +4.21%  "

Set-TextValue $ws.Range("D14") "700.47"
Set-TextValue $ws.Range("E14") "  +3.16%  "

Set-TextValue $ws.Range("D15") "4.182.32"
Set-TextValue $ws.Range("E15") "  +6.43%  "

Set-TextValue $ws.Range("D16") "8.93"
Set-TextValue $ws.Range("E16") "  +3.76%  "

Set-TextValue $ws.Range("D17") "71.997.59"
Set-TextValue $ws.Range("E17") "  +3.75%  "

Set-TextValue $ws.Range("D18") "3.567.74"
Set-TextValue $ws.Range("E18") "  +5.13%  "

Set-TextValue $ws.Range("E19") "  +1.50%  "

Set-TextValue $ws.Range("D20") "18.50"
Set-TextValue $ws.Range("E20") "  +4.70%  "

Set-TextValue $ws.Range("D21") "11.77"
Set-TextValue $ws.Range("E21") "  +4.13%  "

Set-TextValue $ws.Range("D22") "0.931"
Set-TextValue $ws.Range("E22") "  +2.95%  "

Set-TextValue $ws.Range("D23") "5.73"
Set-TextValue $ws.Range("E23") "  +5.77%  "

Set-TextValue $ws.Range("D24") "17.76"
Set-TextValue $ws.Range("E24") "  +3.20%  "

Set-TextValue $ws.Range("D25") "104.55"
Set-TextValue $ws.Range("E25") "  +0.98%  "

Set-TextValue $ws.Range("E26") "  +2.18%  "

Set-TextValue $ws.Range("D27") "2.84"
Set-TextValue $ws.Range("E27") "  +4.11%  "

Set-TextValue $ws.Range("D28") "10.06"
Set-TextValue $ws.Range("E28") "  +4.02%  "

Set-TextValue $ws.Range("D29") "35.29"
Set-TextValue $ws.Range("E29") "  +3.72%  "

Set-TextValue $ws.Range("D30") "9.07"
Set-TextValue $ws.Range("E30") "  +4.08%  "

Set-TextValue $ws.Range("D31") "7.53"
Set-TextValue $ws.Range("E31") "  +8.03%  "

Set-TextValue $ws.Range("E32") "  +15.88%  "

Set-TextValue $ws.Range("D33") "591.63"
Set-TextValue $ws.Range("E33") "  +5.20%  "

Set-TextValue $ws.Range("D34") "11.34"
Set-TextValue $ws.Range("E34") "  +1.79%  "

Set-TextValue $ws.Range("E35") "  +0.98%  "

Set-TextValue $ws.Range("D36") "59.84"
Set-TextValue $ws.Range("E36") "  +2.19%  "

Set-TextValue $ws.Range("E37") "  -0.02%  "

Set-TextValue $ws.Range("D38") "3.683.52"
Set-TextValue $ws.Range("E38") "  -0.07%  "

Set-TextValue $ws.Range("E39") "  +4.34%  "

Set-TextValue $ws.Range("B40") "PEPE"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws.Range("D40") "0.0₃0793"
Set-TextValue $ws.Range("E40") "  +13.47%  "

Set-TextValue $ws.Range("B41") "InjectiveProtocol"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D41") "36.38"
Set-TextValue $ws.Range("E41") "  +2.17%  "

Set-TextValue $ws.Range("D42") "3.50"
Set-TextValue $ws.Range("E42") "  +7.78%  "

Set-TextValue $ws.Range("D43") "2.82"
Set-TextValue $ws.Range("E43") "  +5.30%  "

Set-TextValue $ws.Range("D44") "0.0440"
Set-TextValue $ws.Range("E44") "  +3.97%  "

Set-TextValue $ws.Range("D45") "0.349"
Set-TextValue $ws.Range("E45") "  +2.75%  "

Set-TextValue $ws.Range("D46") "3.42"
Set-TextValue $ws.Range("E46") "  +4.39%  "

Set-TextValue $ws.Range("D47") "2.78"
Set-TextValue $ws.Range("E47") "  +3.76%  "

Set-TextValue $ws.Range("D48") "1.49"
Set-TextValue $ws.Range("E48") "  +5.16%  "

Set-TextValue $ws.Range("D49") "0.133"
Set-TextValue $ws.Range("E49") "  +1.93%  "

Set-TextValue $ws.Range("D50") "0.997"
Set-TextValue $ws.Range("E50") "  -0.54%  "

Set-TextValue $ws.Range("D51") "133.46"
Set-TextValue $ws.Range("E51") "  +0.12%  "
